# This script updates the NATMI ligand-receptor pair output values
# (F2-Gp1ba) to reflect a re-run of the pipeline with updated TPM input
# data, per commit "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.654227
$ws.Range("H2").Value = 4.962681
$ws.Range("I2").Value = 0.4107585939979205
$ws.Range("J2").Value = 0.4107585939979205
$ws.Range("M2").Value = 3.624854
$ws.Range("N2").Value = 10.874562
$ws.Range("O2").Value = 0.2900317783616697
$ws.Range("P2").Value = 0.2900317783616697
$ws.Range("Q2").Value = 5.996331357857999
$ws.Range("R2").Value = 53.966982220722
$ws.Range("S2").Value = 0.119133045494556
$ws.Range("T2").Value = 0.119133045494556

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.654227
$ws.Range("H3").Value = 4.962681
$ws.Range("I3").Value = 0.4107585939979205
$ws.Range("J3").Value = 0.4107585939979205
$ws.Range("O3").Value = 0.3900977855855255
$ws.Range("P3").Value = 0.3900977855855255
$ws.Range("Q3").Value = 8.065169953274999
$ws.Range("R3").Value = 72.586529579475
$ws.Range("S3").Value = 0.1602360179288127
$ws.Range("T3").Value = 0.1602360179288127

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.654227
$ws.Range("H4").Value = 4.962681
$ws.Range("I4").Value = 0.4107585939979205
$ws.Range("J4").Value = 0.4107585939979205
$ws.Range("M4").Value = 2.893069333333333
$ws.Range("N4").Value = 8.679207999999999
$ws.Range("O4").Value = 0.2314802316645793
$ws.Range("P4").Value = 0.2314802316645793
$ws.Range("Q4").Value = 4.785793404071999
$ws.Range("R4").Value = 43.07214063664799
$ws.Range("S4").Value = 0.09508249449685552
$ws.Range("T4").Value = 0.09508249449685552

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.654227
$ws.Range("H5").Value = 4.962681
$ws.Range("I5").Value = 0.4107585939979205
$ws.Range("J5").Value = 0.4107585939979205
$ws.Range("M5").Value = 1.104712
$ws.Range("N5").Value = 3.314136
$ws.Range("O5").Value = 0.08839020438822554
$ws.Range("P5").Value = 0.08839020438822554
$ws.Range("Q5").Value = 1.827444417624
$ws.Range("R5").Value = 16.446999758616
$ws.Range("S5").Value = 0.03630703607769634
$ws.Range("T5").Value = 0.03630703607769634

# Row 6
$ws.Range("I6").Value = 0.3200015957958394
$ws.Range("J6").Value = 0.3200015957958394
$ws.Range("M6").Value = 3.624854
$ws.Range("N6").Value = 10.874562
$ws.Range("O6").Value = 0.2900317783616697
$ws.Range("P6").Value = 0.2900317783616697
$ws.Range("Q6").Value = 4.671443596003999
$ws.Range("R6").Value = 42.042992364036
$ws.Range("S6").Value = 0.09281063190723952
$ws.Range("T6").Value = 0.09281063190723952

# Row 7
$ws.Range("I7").Value = 0.3200015957958394
$ws.Range("J7").Value = 0.3200015957958394
$ws.Range("O7").Value = 0.3900977855855255
$ws.Range("P7").Value = 0.3900977855855255
$ws.Range("S7").Value = 0.1248319139037914
$ws.Range("T7").Value = 0.1248319139037914

# Row 8
$ws.Range("I8").Value = 0.3200015957958394
$ws.Range("J8").Value = 0.3200015957958394
$ws.Range("M8").Value = 2.893069333333333
$ws.Range("N8").Value = 8.679207999999999
$ws.Range("O8").Value = 0.2314802316645793
$ws.Range("P8").Value = 0.2314802316645793
$ws.Range("Q8").Value = 3.728373669669332
$ws.Range("R8").Value = 33.55536302702399
$ws.Range("S8").Value = 0.07407404352785597
$ws.Range("T8").Value = 0.07407404352785597

# Row 9
$ws.Range("I9").Value = 0.3200015957958394
$ws.Range("J9").Value = 0.3200015957958394
$ws.Range("M9").Value = 1.104712
$ws.Range("N9").Value = 3.314136
$ws.Range("O9").Value = 0.08839020438822554
$ws.Range("P9").Value = 0.08839020438822554
$ws.Range("Q9").Value = 1.423671076912
$ws.Range("R9").Value = 12.813039692208
$ws.Range("S9").Value = 0.02828500645695258
$ws.Range("T9").Value = 0.02828500645695258

# Row 10
$ws.Range("G10").Value = 0.8858993333333333
$ws.Range("H10").Value = 2.657698
$ws.Range("I10").Value = 0.2199763179924491
$ws.Range("J10").Value = 0.2199763179924491
$ws.Range("M10").Value = 3.624854
$ws.Range("N10").Value = 10.874562
$ws.Range("O10").Value = 0.2900317783616697
$ws.Range("P10").Value = 0.2900317783616697
$ws.Range("Q10").Value = 3.211255742030666
$ws.Range("R10").Value = 28.901301678276
$ws.Range("S10").Value = 0.06380012270480218
$ws.Range("T10").Value = 0.06380012270480218

# Row 11
$ws.Range("G11").Value = 0.8858993333333333
$ws.Range("H11").Value = 2.657698
$ws.Range("I11").Value = 0.2199763179924491
$ws.Range("J11").Value = 0.2199763179924491
$ws.Range("O11").Value = 0.3900977855855255
$ws.Range("P11").Value = 0.3900977855855255
$ws.Range("Q11").Value = 4.319194817172221
$ws.Range("R11").Value = 38.87275335455
$ws.Range("S11").Value = 0.0858122745301118
$ws.Range("T11").Value = 0.0858122745301118

# Row 12
$ws.Range("G12").Value = 0.8858993333333333
$ws.Range("H12").Value = 2.657698
$ws.Range("I12").Value = 0.2199763179924491
$ws.Range("J12").Value = 0.2199763179924491
$ws.Range("M12").Value = 2.893069333333333
$ws.Range("N12").Value = 8.679207999999999
$ws.Range("O12").Value = 0.2314802316645793
$ws.Range("P12").Value = 0.2314802316645793
$ws.Range("Q12").Value = 2.56296819368711
$ws.Range("R12").Value = 23.066713743184
$ws.Range("S12").Value = 0.05092016904961329
$ws.Range("T12").Value = 0.05092016904961329

# Row 13
$ws.Range("G13").Value = 0.8858993333333333
$ws.Range("H13").Value = 2.657698
$ws.Range("I13").Value = 0.2199763179924491
$ws.Range("J13").Value = 0.2199763179924491
$ws.Range("M13").Value = 1.104712
$ws.Range("N13").Value = 3.314136
$ws.Range("O13").Value = 0.08839020438822554
$ws.Range("P13").Value = 0.08839020438822554
$ws.Range("Q13").Value = 0.9786636243253332
$ws.Range("R13").Value = 8.807972618928
$ws.Range("S13").Value = 0.01944375170792187
$ws.Range("T13").Value = 0.01944375170792187

# Row 14
$ws.Range("G14").Value = 0.1983963333333333
$ws.Range("H14").Value = 0.595189
$ws.Range("I14").Value = 0.04926349221379096
$ws.Range("J14").Value = 0.04926349221379096
$ws.Range("M14").Value = 3.624854
$ws.Range("N14").Value = 10.874562
$ws.Range("O14").Value = 0.2900317783616697
$ws.Range("P14").Value = 0.2900317783616697
$ws.Range("Q14").Value = 0.7191577424686665
$ws.Range("R14").Value = 6.472419682217999
$ws.Range("S14").Value = 0.01428797825507206
$ws.Range("T14").Value = 0.01428797825507206

# Row 15
$ws.Range("G15").Value = 0.1983963333333333
$ws.Range("H15").Value = 0.595189
$ws.Range("I15").Value = 0.04926349221379096
$ws.Range("J15").Value = 0.04926349221379096
$ws.Range("O15").Value = 0.3900977855855255
$ws.Range("P15").Value = 0.3900977855855255
$ws.Range("Q15").Value = 0.9672796698638887
$ws.Range("R15").Value = 8.705517028774999
$ws.Range("S15").Value = 0.01921757922280963
$ws.Range("T15").Value = 0.01921757922280963

# Row 16
$ws.Range("G16").Value = 0.1983963333333333
$ws.Range("H16").Value = 0.595189
$ws.Range("I16").Value = 0.04926349221379096
$ws.Range("J16").Value = 0.04926349221379096
$ws.Range("M16").Value = 2.893069333333333
$ws.Range("N16").Value = 8.679207999999999
$ws.Range("O16").Value = 0.2314802316645793
$ws.Range("P16").Value = 0.2314802316645793
$ws.Range("Q16").Value = 0.5739743478124444
$ws.Range("R16").Value = 5.165769130311999
$ws.Range("S16").Value = 0.01140352459025453
$ws.Range("T16").Value = 0.01140352459025453

# Row 17
$ws.Range("G17").Value = 0.1983963333333333
$ws.Range("H17").Value = 0.595189
$ws.Range("I17").Value = 0.04926349221379096
$ws.Range("J17").Value = 0.04926349221379096
$ws.Range("M17").Value = 1.104712
$ws.Range("N17").Value = 3.314136
$ws.Range("O17").Value = 0.08839020438822554
$ws.Range("P17").Value = 0.08839020438822554
$ws.Range("Q17").Value = 0.2191708101893333
$ws.Range("R17").Value = 1.972537291704
$ws.Range("S17").Value = 0.00435441014565474
$ws.Range("T17").Value = 0.00435441014565474

